# Applies the "added schema files from jenkins job" commit:
#  1. On the "Common" sheet (index 2), inserts a new row at position 70
#     ("VSTAT License File") shifting all subsequent rows (merged cells,
#     data validations, comments, ...) down by one.
#  2. On the "Credentials" sheet (index 7), rewords several cell comments
#     and the text of cell A57.

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# 1. "Common" sheet - insert "VSTAT License File" row before row 70
# -----------------------------------------------------------------------
$ws = $wb.Worksheets.Item(2)

$lastRow = 129

# Preserve the comment text currently attached to rows 70..lastRow (col A)
# so that we can re-apply it, shifted down by one row, after the insert.
# (A native Rows.Insert() moves cell values/formatting/validation, but
# comments stay pinned to their original row, so we move them ourselves.)
$savedComments = @{}
for ($r = 70; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Comment -ne $null) {
        $savedComments[$r] = $cell.Comment.Text()
    }
}

# Insert the new row; this shifts cell values, styles, merged cells and
# data validations down by one automatically.
$ws.Rows.Item(70).Insert()

# The inserted row is blank/default-formatted; copy the formatting (only)
# from the row directly below (the old row 70, now row 71) so the new row
# matches the look of its neighbours.
$ws.Cells.Item(71, 1).Resize(1, 2).Copy()
$ws.Cells.Item(70, 1).Resize(1, 2).PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(70, 1).Value = "VSTAT License File"

# Re-apply the saved comments, shifted down by one row.
for ($r = $lastRow; $r -ge 70; $r--) {
    $destRow = $r + 1
    $destCell = $ws.Cells.Item($destRow, 1)
    if ($savedComments.ContainsKey($r)) {
        if ($destCell.Comment -ne $null) {
            $destCell.Comment.Text($savedComments[$r]) | Out-Null
        } else {
            $destCell.AddComment($savedComments[$r]) | Out-Null
        }
    } elseif ($destCell.Comment -ne $null) {
        $destCell.Comment.Delete()
    }
}

# Finally, give the new row its own comment.
$newCell = $ws.Cells.Item(70, 1)
if ($newCell.Comment -ne $null) {
    $newCell.Comment.Delete()
}
$newCell.AddComment("Optional License File for Elasticsearch [default: ]") | Out-Null

# -----------------------------------------------------------------------
# 2. "Credentials" sheet - reword several comments + cell A57
# -----------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item(7)

$ws7.Cells.Item(57, 1).Value = "NETCONF Manager VM password for running sudo commands, and will be used for the installation of NETCONF Manager."

function Set-CommentText($sheet, $addr, $text) {
    $cell = $sheet.Range($addr)
    if ($cell.Comment -ne $null) {
        $cell.Comment.Text($text) | Out-Null
    } else {
        $cell.AddComment($text) | Out-Null
    }
}

Set-CommentText $ws7 "A10" "VSD Username will be used for logging into VSD command line. Used for both Install and Upgrade procedures. [default: root]"
Set-CommentText $ws7 "A11" "VSD password will be used for logging into the command line. Used for both Install and Upgrade procedures. [default: Alcateldc]"
Set-CommentText $ws7 "A12" "VSC Username will be used for logging into command line (should have admin privileges). Used for upgrade procedure only [default: ]"
Set-CommentText $ws7 "A13" "VSC password will be used for logging into the command line. Used for upgrade procedure only [default: ]"
Set-CommentText $ws7 "A15" "ElasticSearch (Stats) Username will be used for logging into command line. Used for both Install and Upgrade procedures. [default: ]"
Set-CommentText $ws7 "A16" "ElasticSearch (Stats) password will be used for logging into the command line. Used for both Install and Upgrade procedures. [default: ]"
Set-CommentText $ws7 "A17" "ElasticSearch (Stats) root password required for VSTAT Upgrade only [default: ]"
Set-CommentText $ws7 "A19" "This VSD Username(also known as csproot user). Used for both Install and Upgrade procedures. Must have csproot privileges. [default: csproot]"
Set-CommentText $ws7 "A20" "This VSD password(also known as csproot password) will be used for API authentication. Used for both Install and Upgrade procedures. Must have csproot privileges. [default: csproot]"
Set-CommentText $ws7 "A21" "This VSD Mysql password. Used for both Install and Upgrade procedures. [default: ]"
Set-CommentText $ws7 "A37" "Username for OpenStack. [default: ]"
Set-CommentText $ws7 "A40" "vCenter Username. [default: ]"
Set-CommentText $ws7 "A43" "Username for Compute node to install VRS. [default: root]"
Set-CommentText $ws7 "A44" "Password for Compute node, and will be used for installation of VRS [default: ]"
Set-CommentText $ws7 "A54" "NFS username to login into command line, and will be used for NFS configuration. Default user is root. [default: root]"
Set-CommentText $ws7 "A56" "Username for NETCONF Manager VM, and will be used for the installation of NETCONF Manager. Default user is root. [default: root]"
Set-CommentText $ws7 "A58" "Username for NETCONF Manager user, and will be used for the installation of NETCONF Manager. [default: netconfmgr]"
Set-CommentText $ws7 "A59" "Password for NETCONF manager user, and will be used for the installation of NETCONF Manager. [default: password]"
Set-CommentText $ws7 "A61" "Username for SMTP Server, and will be used for Email health report."
Set-CommentText $ws7 "A62" "Password for SMTP Server, and will be used for Email health report."
Set-CommentText $ws7 "A64" "Username for the monit mail server."
Set-CommentText $ws7 "A67" "Username for NUH notification application, and will be used for installation of NUH."
Set-CommentText $ws7 "A68" "Password for NUH notification application, and will be used for installation of NUH."
Set-CommentText $ws7 "A69" "Username for NUH notification application, and will be used for installation of NUH."
Set-CommentText $ws7 "A70" "Password for NUH notification application, and will be used for installation of NUH."
